$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $text) {
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

Set-TextValue $ws.Cells(2, 4) "42.482.63"
$ws.Range("E2").Value = "  +0.51%  "
Set-TextValue $ws.Cells(3, 4) "2.277.99"
$ws.Range("E3").Value = "  -0.72%  "
Set-TextValue $ws.Cells(4, 4) "1.02"
$ws.Range("E4").Value = "  +1.82%  "
Set-TextValue $ws.Cells(5, 4) "311.71"
$ws.Range("E5").Value = "  -1.61%  "
Set-TextValue $ws.Cells(6, 4) "101.08"
$ws.Range("E6").Value = "  -1.42%  "
Set-TextValue $ws.Cells(7, 4) "0.620"
$ws.Range("E7").Value = "  -0.38%  "
Set-TextValue $ws.Cells(8, 4) "1.00"
Set-TextValue $ws.Cells(9, 4) "0.593"
$ws.Range("E9").Value = "  -2.21%  "
Set-TextValue $ws.Cells(10, 4) "38.44"
$ws.Range("E10").Value = "  -2.77%  "
Set-TextValue $ws.Cells(11, 4) "0.0894"
$ws.Range("E11").Value = "  -1.31%  "
Set-TextValue $ws.Cells(12, 4) "8.19"
$ws.Range("E12").Value = "  -2.36%  "
$ws.Range("E13").Value = "  +1.71%  "
$ws.Range("E14").Value = "  +2.01%  "
Set-TextValue $ws.Cells(15, 4) "14.96"
$ws.Range("E15").Value = "  -1.80%  "
Set-TextValue $ws.Cells(16, 4) "2.625.54"
$ws.Range("E16").Value = "  -0.63%  "
Set-TextValue $ws.Cells(17, 4) "2.320.02"
$ws.Range("E17").Value = "  +1.56%  "
Set-TextValue $ws.Cells(18, 4) "42.360.35"
$ws.Range("E18").Value = "  -0.04%  "
Set-TextValue $ws.Cells(19, 4) "7.21"
$ws.Range("E19").Value = "  -2.57%  "
$ws.Range("E20").Value = "  -1.41%  "
Set-TextValue $ws.Cells(21, 4) "13.33"
$ws.Range("E21").Value = "  +7.83%  "
Set-TextValue $ws.Cells(22, 4) "72.70"
$ws.Range("E22").Value = "  -0.96%  "
Set-TextValue $ws.Cells(23, 4) "3.47"
$ws.Range("E23").Value = "  -2.28%  "
Set-TextValue $ws.Cells(24, 4) "261.21"
$ws.Range("E24").Value = "  -5.34%  "
Set-TextValue $ws.Cells(25, 4) "2.15"
$ws.Range("E25").Value = "  -4.77%  "
Set-TextValue $ws.Cells(26, 4) "1.00"
$ws.Range("E26").Value = "  +0.12%  "
Set-TextValue $ws.Cells(27, 4) "10.59"
$ws.Range("E27").Value = "  -2.08%  "
$ws.Range("E28").Value = "  -1.34%  "
Set-TextValue $ws.Cells(29, 4) "6.79"
$ws.Range("E29").Value = "  +12.69%  "
Set-TextValue $ws.Cells(30, 4) "22.19"
$ws.Range("E30").Value = "  -2.45%  "
Set-TextValue $ws.Cells(31, 4) "35.64"
$ws.Range("E31").Value = "  -4.74%  "
Set-TextValue $ws.Cells(32, 4) "164.63"
$ws.Range("E32").Value = "  -0.98%  "
Set-TextValue $ws.Cells(33, 4) "0.0857"
$ws.Range("E33").Value = "  -1.94%  "
Set-TextValue $ws.Cells(34, 4) "0.129"
$ws.Range("E34").Value = "  -3.29%  "
$ws.Range("E35").Value = "  -1.52%  "
$ws.Range("E36").Value = "  -5.45%  "
Set-TextValue $ws.Cells(37, 4) "4.45"
$ws.Range("E37").Value = "  -2.93%  "
$ws.Range("E38").Value = "  -4.55%  "
Set-TextValue $ws.Cells(39, 4) "3.67"
$ws.Range("E39").Value = "  -0.66%  "
$ws.Range("E40").Value = "  -3.76%  "
Set-TextValue $ws.Cells(41, 4) "1.56"
$ws.Range("E41").Value = "  +4.61%  "
Set-TextValue $ws.Cells(42, 4) "95.97"
$ws.Range("E42").Value = "  -0.54%  "
Set-TextValue $ws.Cells(45, 4) "0.224"
$ws.Range("E45").Value = "  -0.40%  "
Set-TextValue $ws.Cells(46, 4) "11.85"
$ws.Range("E46").Value = "  -0.93%  "
Set-TextValue $ws.Cells(47, 4) "1.704.00"
$ws.Range("E47").Value = "  +6.62%  "
Set-TextValue $ws.Cells(48, 4) "78.84"
$ws.Range("E48").Value = "  +0.15%  "
Set-TextValue $ws.Cells(49, 4) "109.80"
$ws.Range("E49").Value = "  -2.59%  "
Set-TextValue $ws.Cells(50, 4) "8.67"
$ws.Range("E50").Value = "  -3.28%  "
Set-TextValue $ws.Cells(51, 4) "5.14"
$ws.Range("E51").Value = "  -2.44%  "

# Rows 43 and 44 swap coin identity (MultiversX <-> FirstDigitalUSD) and receive updated price/volume
$ws.Range("B43").Value = "FirstDigitalUSD"
$ws.Range("C43").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
Set-TextValue $ws.Cells(43, 4) "1.01"
$ws.Range("E43").Value = "  +0.19%  "

$ws.Range("B44").Value = "MultiversX"
$ws.Range("C44").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
Set-TextValue $ws.Cells(44, 4) "68.57"
$ws.Range("E44").Value = "  -1.72%  "

